$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.32689766666667
$ws.Range("H2").Value = 30.980693
$ws.Range("I2").Value = 0.2044815006034941
$ws.Range("J2").Value = 0.204481500603494
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 1161.914285313552
$ws.Range("R2").Value = 10457.22856782197
$ws.Range("S2").Value = 0.0669740912558422
$ws.Range("T2").Value = 0.06697409125584218
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.32689766666667
$ws.Range("H3").Value = 30.980693
$ws.Range("I3").Value = 0.2044815006034941
$ws.Range("J3").Value = 0.204481500603494
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 1097.898610868313
$ws.Range("R3").Value = 9881.087497814813
$ws.Range("S3").Value = 0.06328415330061449
$ws.Range("T3").Value = 0.06328415330061449
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.32689766666667
$ws.Range("H4").Value = 30.980693
$ws.Range("I4").Value = 0.2044815006034941
$ws.Range("J4").Value = 0.204481500603494
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 1287.677964514604
$ws.Range("R4").Value = 11589.10168063144
$ws.Range("S4").Value = 0.07422325604703739
$ws.Range("T4").Value = 0.07422325604703739
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.37031933333333
$ws.Range("H5").Value = 88.110958
$ws.Range("I5").Value = 0.5815577111671272
$ws.Range("J5").Value = 0.5815577111671272
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 3304.554252316512
$ws.Range("R5").Value = 29740.9882708486
$ws.Range("S5").Value = 0.190478351847445
$ws.Range("T5").Value = 0.190478351847445
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.37031933333333
$ws.Range("H6").Value = 88.110958
$ws.Range("I6").Value = 0.5815577111671272
$ws.Range("J6").Value = 0.5815577111671272
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 3122.489816172809
$ws.Range("R6").Value = 28102.40834555528
$ws.Range("S6").Value = 0.1799839459219328
$ws.Range("T6").Value = 0.1799839459219329
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.37031933333333
$ws.Range("H7").Value = 88.110958
$ws.Range("I7").Value = 0.5815577111671272
$ws.Range("J7").Value = 0.5815577111671272
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 3662.233735341937
$ws.Range("R7").Value = 32960.10361807743
$ws.Range("S7").Value = 0.2110954133977493
$ws.Range("T7").Value = 0.2110954133977493
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.80562866666667
$ws.Range("H8").Value = 32.416886
$ws.Range("I8").Value = 0.2139607882293788
$ws.Range("J8").Value = 0.2139607882293788
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 1215.777933979104
$ws.Range("R8").Value = 10942.00140581193
$ws.Range("S8").Value = 0.07007885463356914
$ws.Range("T8").Value = 0.07007885463356914
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.80562866666667
$ws.Range("H9").Value = 32.416886
$ws.Range("I9").Value = 0.2139607882293788
$ws.Range("J9").Value = 0.2139607882293788
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 1148.794641490958
$ws.Range("R9").Value = 10339.15177341863
$ws.Range("S9").Value = 0.06621785972161899
$ws.Range("T9").Value = 0.06621785972161899
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.80562866666667
$ws.Range("H10").Value = 32.416886
$ws.Range("I10").Value = 0.2139607882293788
$ws.Range("J10").Value = 0.2139607882293788
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 1347.371725363986
$ws.Range("R10").Value = 12126.34552827588
$ws.Range("S10").Value = 0.07766407387419066
$ws.Range("T10").Value = 0.07766407387419066
